$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.993232250213623
$ws.Range("B1").Value = 1.986650466918945
$ws.Range("C1").Value = 1.905375957489014
$ws.Range("D1").Value = 1.753554224967957
$ws.Range("E1").Value = 1.534249663352966
